$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the city names in parentheses from the state names in column A
$ws.Range("A2").Value = "South Australia"
$ws.Range("A3").Value = "Queensland"
$ws.Range("A4").Value = "Australian Capital Territory"
$ws.Range("A5").Value = "Northern Territory"
$ws.Range("A6").Value = "Tasmania"
$ws.Range("A7").Value = "Victoria "
$ws.Range("A8").Value = "Western Australia"
$ws.Range("A9").Value = "New South Wales"

# Update the selected cell to match the saved selection state
$ws.Range("B3").Select()
